$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 176040
$ws.Range("C4").Value = 166003
$ws.Range("C5").Value = 10037
$ws.Range("C8").Value = 64.72
